# Regenerate merged AHB files
#
# 1) Rename the "_old" / "_new" header-name suffixes to the new
#    release tags "_FV2304" / "_FV2310" (row 1 header cells only -
#    these substrings only occur in the header row of the sheet).
# 2) Turn the used range A1:U63 into a native Excel Table (ListObject).
# 3) Freeze the header row (pane split under row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header labels -------------------------------------------------
$used = $ws.UsedRange
$used.Replace("_old", "_FV2304", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$used.Replace("_new", "_FV2310", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)

# --- 2. Convert the range into a table ---------------------------------------
$dataRange = $ws.Range("A1:U63")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# --- 3. Freeze the top (header) row -------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
